$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Stage a full snapshot (value + number/font/fill formatting) of the current
# A1:G1 contents in a scratch row (row 50) so that the cyclic rearrangement
# of columns A-H in row 1 can be performed without any cell clobbering its
# own, not-yet-read, source value. H1 is a formula, handled separately below.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy($ws.Range("A50"))
$ws.Range("B1").Copy($ws.Range("B50"))
$ws.Range("C1").Copy($ws.Range("C50"))
$ws.Range("D1").Copy($ws.Range("D50"))
$ws.Range("E1").Copy($ws.Range("E50"))
$ws.Range("F1").Copy($ws.Range("F50"))
$ws.Range("G1").Copy($ws.Range("G50"))

# ---------------------------------------------------------------------------
# Write the new arrangement into row 1 from the staged snapshot.
#   A1 <- old F1 (empty, red-font style)
#   B1 <- old G1 (empty, red-font+fill style)
#   C1 <- old A1 ("String")
#   D1 <- formula, now referencing the relocated C1 (no style, like old H1)
#   E1 <- old B1 (1)
#   F1 <- old C1 (1.1)
#   G1 <- old D1 (date, styled)
#   H1 <- old E1 (currency, styled)
# ---------------------------------------------------------------------------
$ws.Range("F50").Copy($ws.Range("A1"))
$ws.Range("A1").ClearContents()
$ws.Range("G50").Copy($ws.Range("B1"))
$ws.Range("B1").ClearContents()
$ws.Range("A50").Copy($ws.Range("C1"))
$ws.Range("D1").ClearFormats()
$ws.Range("D1").Formula = "=CONCAT(C1,C1)"
$ws.Range("B50").Copy($ws.Range("E1"))
$ws.Range("C50").Copy($ws.Range("F1"))
$ws.Range("D50").Copy($ws.Range("G1"))
$ws.Range("E50").Copy($ws.Range("H1"))

# Formulas in I1/J1 must now reference the relocated operand cells.
$ws.Range("I1").Formula = "=E1+E1"
$ws.Range("J1").Formula = "=F1+F1"

# ---------------------------------------------------------------------------
# Clean up the scratch staging row.
# ---------------------------------------------------------------------------
$ws.Range("A50:G50").Clear()
